$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the caption/footnote text that was previously missing for the image in A5
$ws.Range("A5").Value = "Figura 2: Struttura del frame Modbus in modalità RTU, consultato il 10 settembre 2023,  https://development.libelium.com/modbus_networking_guide/introduction"

# Move the active selection to A7 (was B7)
$ws.Range("A7").Select()
